$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 25.92718227809252
$ws.Range("B1").Value = 13.65578108662694
$ws.Range("C1").Value = 27.71006605668829
$ws.Range("D1").Value = 7.910039432272094
$ws.Range("E1").Value = 26.7473686363864
$ws.Range("F1").Value = 25.94674647443332
$ws.Range("G1").Value = 21.02766578514784
$ws.Range("H1").Value = 9.006345202225754
$ws.Range("I1").Value = 10.30789802559425
$ws.Range("J1").Value = 26.4828424409353
$ws.Range("K1").Value = 3.071574792433862
$ws.Range("L1").Value = 8.735876604905631
$ws.Range("M1").Value = 29.533812194501
$ws.Range("N1").Value = 1.727206190720807
$ws.Range("O1").Value = 2.32165797744991
$ws.Range("A2").Value = 20.15770125058977
$ws.Range("B2").Value = 1.450994709476593
$ws.Range("C2").Value = 8.784064384793648
$ws.Range("D2").Value = 11.3309660173357
$ws.Range("E2").Value = 26.04615144564942
$ws.Range("F2").Value = 27.62830356771391
$ws.Range("G2").Value = 19.79424949414
$ws.Range("H2").Value = 29.20483818295027
$ws.Range("I2").Value = 4.913380421086946
$ws.Range("J2").Value = 26.28852095168888
$ws.Range("K2").Value = 17.26487618154057
$ws.Range("L2").Value = 1.218753656114688
$ws.Range("M2").Value = 18.22760923647114
$ws.Range("N2").Value = 24.8416834608759
$ws.Range("O2").Value = 6.734479217225015
$ws.Range("A3").Value = 20.86540690842241
$ws.Range("B3").Value = 3.450433580973494
$ws.Range("C3").Value = 7.355774208074092
$ws.Range("D3").Value = 27.27335670502888
$ws.Range("E3").Value = 18.23571952212412
$ws.Range("F3").Value = 19.60901616593736
$ws.Range("G3").Value = 28.94050355056864
$ws.Range("H3").Value = 18.67931166116659
$ws.Range("I3").Value = 2.481793334315846
$ws.Range("J3").Value = 23.08853847460831
$ws.Range("K3").Value = 11.13888500375606
$ws.Range("L3").Value = 26.44479386430688
$ws.Range("M3").Value = 26.51360782181713
$ws.Range("N3").Value = 22.03459177630727
$ws.Range("O3").Value = 20.99471556347149
$ws.Range("A4").Value = 7.373482331144545
$ws.Range("B4").Value = 20.66575447224554
$ws.Range("C4").Value = 11.29363265627161
$ws.Range("D4").Value = 22.21255412074022
$ws.Range("E4").Value = 26.76403648387036
$ws.Range("F4").Value = 13.46733099186202
$ws.Range("G4").Value = 29.8646498160022
$ws.Range("H4").Value = 26.65857581369774
$ws.Range("I4").Value = 8.138160518515273
$ws.Range("J4").Value = 24.54383608154318
$ws.Range("K4").Value = 1.616779998722623
$ws.Range("L4").Value = 25.57788298418383
$ws.Range("M4").Value = 29.84387676462801
$ws.Range("N4").Value = 21.24367316605093
$ws.Range("O4").Value = 6.358287878375246
$ws.Range("A5").Value = 14.4889673995826
$ws.Range("B5").Value = 20.44250054911469
$ws.Range("C5").Value = 4.191709995990229
$ws.Range("D5").Value = 11.9164784218025
$ws.Range("E5").Value = 4.459209369497922
$ws.Range("F5").Value = 14.24681002363393
$ws.Range("G5").Value = 12.42363938345982
$ws.Range("H5").Value = 2.737444558477874
$ws.Range("I5").Value = 29.35838722390753
$ws.Range("J5").Value = 5.19944245095222
$ws.Range("K5").Value = 28.66840657661794
$ws.Range("L5").Value = 29.52545579582688
$ws.Range("M5").Value = 22.41394137553782
$ws.Range("N5").Value = 12.88352701473961
$ws.Range("O5").Value = 14.89212698844447
$ws.Range("A6").Value = 12.58851111904495
$ws.Range("B6").Value = 17.25593964634868
$ws.Range("C6").Value = 27.08093052878391
$ws.Range("D6").Value = 17.10026196179122
$ws.Range("E6").Value = 10.49755827512746
$ws.Range("F6").Value = 2.051523546837716
$ws.Range("G6").Value = 24.17875529724298
$ws.Range("H6").Value = 22.10064362575594
$ws.Range("I6").Value = 3.590911208515566
$ws.Range("J6").Value = 26.44669027346433
$ws.Range("K6").Value = 21.60741330156816
$ws.Range("L6").Value = 22.41943350393503
$ws.Range("M6").Value = 16.60143450419001
$ws.Range("N6").Value = 5.665820019811048
$ws.Range("O6").Value = 25.2866829055219
$ws.Range("A7").Value = 12.08259603383189
$ws.Range("B7").Value = 15.23282665543521
$ws.Range("C7").Value = 8.429325914113431
$ws.Range("D7").Value = 11.54686572319147
$ws.Range("E7").Value = 16.82423756822132
$ws.Range("F7").Value = 16.66775687050954
$ws.Range("G7").Value = 10.54164788901882
$ws.Range("H7").Value = 5.937007833692081
$ws.Range("I7").Value = 9.850296071111124
$ws.Range("J7").Value = 12.32406277207107
$ws.Range("K7").Value = 6.495115581432477
$ws.Range("L7").Value = 13.16966543319803
$ws.Range("M7").Value = 11.16623041499016
$ws.Range("N7").Value = 19.33771617287197
$ws.Range("O7").Value = 3.703895796245941
$ws.Range("A8").Value = 15.243974615202
$ws.Range("B8").Value = 28.15791846773456
$ws.Range("C8").Value = 21.28625891985301
$ws.Range("D8").Value = 12.90898736051674
$ws.Range("E8").Value = 3.558429584525522
$ws.Range("F8").Value = 1.222139359567321
$ws.Range("G8").Value = 10.29681449979964
$ws.Range("H8").Value = 24.48553650820482
$ws.Range("I8").Value = 11.6492200706358
$ws.Range("J8").Value = 27.80417192144546
$ws.Range("K8").Value = 20.04650391289548
$ws.Range("L8").Value = 4.410772514586058
$ws.Range("M8").Value = 25.30041141519226
$ws.Range("N8").Value = 9.262166491454508
$ws.Range("O8").Value = 1.489841307158386
$ws.Range("A9").Value = 3.867765299773526
$ws.Range("B9").Value = 4.957555062332055
$ws.Range("C9").Value = 8.180042760935141
$ws.Range("D9").Value = 22.73067325031145
$ws.Range("E9").Value = 26.26517405928889
$ws.Range("F9").Value = 2.096569759411483
$ws.Range("G9").Value = 5.63012458336043
$ws.Range("H9").Value = 10.65610111097597
$ws.Range("I9").Value = 9.023389875049221
$ws.Range("J9").Value = 20.0727060877417
$ws.Range("K9").Value = 26.68963346623893
$ws.Range("L9").Value = 16.19644497369756
$ws.Range("M9").Value = 6.284949053738711
$ws.Range("N9").Value = 17.2438913143201
$ws.Range("O9").Value = 21.50448201789163
$ws.Range("A10").Value = 12.78931290264147
$ws.Range("B10").Value = 27.95513288348207
$ws.Range("C10").Value = 5.542397677384428
$ws.Range("D10").Value = 10.04420863906235
$ws.Range("E10").Value = 17.77647911017322
$ws.Range("F10").Value = 9.602590844918373
$ws.Range("G10").Value = 13.35001869729106
$ws.Range("H10").Value = 1.722129677228411
$ws.Range("I10").Value = 11.03198188628979
$ws.Range("J10").Value = 16.76329562101495
$ws.Range("K10").Value = 24.01852074106805
$ws.Range("L10").Value = 22.39552612028641
$ws.Range("M10").Value = 10.51064409650425
$ws.Range("N10").Value = 14.34315099529383
$ws.Range("O10").Value = 8.901257387151126
$ws.Range("A11").Value = 13.14758647698259
$ws.Range("B11").Value = 20.48030952289237
$ws.Range("C11").Value = 6.078492349624941
$ws.Range("D11").Value = 15.66942298570057
$ws.Range("E11").Value = 5.861167315695384
$ws.Range("F11").Value = 22.92140451517371
$ws.Range("G11").Value = 24.77410383624736
$ws.Range("H11").Value = 22.93041402784607
$ws.Range("I11").Value = 23.86734851940706
$ws.Range("J11").Value = 8.419760211949587
$ws.Range("K11").Value = 27.91046424709499
$ws.Range("L11").Value = 14.28698928920588
$ws.Range("M11").Value = 27.60853096486451
$ws.Range("N11").Value = 21.05788838293685
$ws.Range("O11").Value = 29.9214720156209
$ws.Range("A12").Value = 12.87284893647637
$ws.Range("B12").Value = 7.321300924394389
$ws.Range("C12").Value = 20.37823911267997
$ws.Range("D12").Value = 7.931702599577209
$ws.Range("E12").Value = 15.07433981712331
$ws.Range("F12").Value = 22.51562092900768
$ws.Range("G12").Value = 19.68029726492231
$ws.Range("H12").Value = 2.380982017536825
$ws.Range("I12").Value = 5.723296727042653
$ws.Range("J12").Value = 28.91971976610722
$ws.Range("K12").Value = 7.795921389559149
$ws.Range("L12").Value = 21.368981248342
$ws.Range("M12").Value = 16.66751768394892
$ws.Range("N12").Value = 5.191305531889839
$ws.Range("O12").Value = 16.30986653827399
$ws.Range("A13").Value = 11.38409583899061
$ws.Range("B13").Value = 4.393729294287832
$ws.Range("C13").Value = 6.661436790191807
$ws.Range("D13").Value = 4.282348756470161
$ws.Range("E13").Value = 28.09506930151939
$ws.Range("F13").Value = 23.20557937634998
$ws.Range("G13").Value = 4.130157728067699
$ws.Range("H13").Value = 2.586842153440894
$ws.Range("I13").Value = 15.23422438074079
$ws.Range("J13").Value = 26.06757407379948
$ws.Range("K13").Value = 15.81327078037232
$ws.Range("L13").Value = 27.77464927717223
$ws.Range("M13").Value = 25.78421838537085
$ws.Range("N13").Value = 4.011239984298193
$ws.Range("O13").Value = 20.28724031974253
$ws.Range("A14").Value = 6.820272518211079
$ws.Range("B14").Value = 26.15736283146289
$ws.Range("C14").Value = 14.44350209485315
$ws.Range("D14").Value = 20.23658713028516
$ws.Range("E14").Value = 21.90976619391834
$ws.Range("F14").Value = 20.72848403701386
$ws.Range("G14").Value = 17.3394669612568
$ws.Range("H14").Value = 28.3753912516888
$ws.Range("I14").Value = 27.47846726079839
$ws.Range("J14").Value = 21.72876512411175
$ws.Range("K14").Value = 16.6092475771182
$ws.Range("L14").Value = 5.505391559928768
$ws.Range("M14").Value = 4.306570994143156
$ws.Range("N14").Value = 24.67444124306629
$ws.Range("O14").Value = 19.18843641230858
$ws.Range("A15").Value = 10.61491848661818
$ws.Range("B15").Value = 1.301329795507729
$ws.Range("C15").Value = 23.17470076406565
$ws.Range("D15").Value = 1.865508259636456
$ws.Range("E15").Value = 25.28078642068363
$ws.Range("F15").Value = 8.300469936984037
$ws.Range("G15").Value = 23.97458171890923
$ws.Range("H15").Value = 8.878301832807477
$ws.Range("I15").Value = 12.88898550789792
$ws.Range("J15").Value = 29.21181414349279
$ws.Range("K15").Value = 21.08774988688578
$ws.Range("L15").Value = 15.5330257791087
$ws.Range("M15").Value = 14.43782707359506
$ws.Range("N15").Value = 17.54777872928039
$ws.Range("O15").Value = 16.72090132682687
